$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "Changed the layout of the analytics page and added some more data to the page"
$row1.Cells.Item(2).Range.Text = "15/01/2025"

$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Added a show/hide password button on login, register, confirm password, and update profile pages"
$row2.Cells.Item(2).Range.Text = "16/01/2025"
